$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9444
$ws.Range("D2").Value = 8346
$ws.Range("E2").Value = 0.8837357052096569
$ws.Range("F2").Value = 0.8805655201519308
$ws.Range("G2").Value = 0.0970076441535767
$ws.Range("H2").Value = 0.08542158663280766
$ws.Range("I2").Value = 41076824.4099903
$ws.Range("J2").Value = 14343580.15165115
$ws.Range("L2").Value = 14343580.15165115
$ws.Range("M2").Value = 55420404.56164145
$ws.Range("N2").Value = 798913812.9572
$ws.Range("O2").Value = 781214005.9532001
$ws.Range("P2").Value = 0.01795385174097568
$ws.Range("Q2").Value = 0.01836062851196554

$ws.Range("C3").Value = 9640
$ws.Range("D3").Value = 8568
$ws.Range("E3").Value = 0.8887966804979253
$ws.Range("F3").Value = 0.8864059590316573
$ws.Range("G3").Value = 0.09542172038036506
$ws.Range("H3").Value = 0.08458238156620813
$ws.Range("I3").Value = 43142786.81432747
$ws.Range("J3").Value = 15106644.34690293
$ws.Range("L3").Value = 15106644.34690293
$ws.Range("M3").Value = 58249431.1612304
$ws.Range("N3").Value = 837045574.001528
$ws.Range("O3").Value = 819565397.977458
$ws.Range("P3").Value = 0.01804757687766635
$ws.Range("Q3").Value = 0.01843250628221183

$ws.Range("C4").Value = 9845
$ws.Range("D4").Value = 8756
$ws.Range("E4").Value = 0.8893854748603351
$ws.Range("F4").Value = 0.888212619192534
$ws.Range("G4").Value = 0.09412919338671218
$ws.Range("H4").Value = 0.08360673740049218
$ws.Range("I4").Value = 45118047.21266638
$ws.Range("J4").Value = 15774352.41540857
$ws.Range("L4").Value = 15774352.41540857
$ws.Range("M4").Value = 60892399.62807495
$ws.Range("N4").Value = 874986044.7346259
$ws.Range("O4").Value = 857537096.7286721
$ws.Range("P4").Value = 0.01802811886010453
$ws.Range("Q4").Value = 0.0183949504640493

$ws.Range("C5").Value = 10031
$ws.Range("D5").Value = 8913
$ws.Range("E5").Value = 0.8885455089223407
$ws.Range("F5").Value = 0.8865128307141437
$ws.Range("G5").Value = 0.09321308070362111
$ws.Range("H5").Value = 0.08263459203415308
$ws.Range("I5").Value = 47213184.81612386
$ws.Range("J5").Value = 16488847.50894834
$ws.Range("L5").Value = 16488847.50894834
$ws.Range("M5").Value = 63702032.3250722
$ws.Range("N5").Value = 914432692.5229203
$ws.Range("O5").Value = 896946587.059947
$ws.Range("P5").Value = 0.01803177822028169
$ws.Range("Q5").Value = 0.01838331038528866

$ws.Range("C6").Value = 10225
$ws.Range("D6").Value = 9080
$ws.Range("E6").Value = 0.8880195599022005
$ws.Range("F6").Value = 0.8855080944021845
$ws.Range("G6").Value = 0.09213732329015105
$ws.Range("H6").Value = 0.08158834556997968
$ws.Range("I6").Value = 49406731.03858929
$ws.Range("J6").Value = 17214575.96035406
$ws.Range("L6").Value = 17214575.96035406
$ws.Range("M6").Value = 66621306.99894334
$ws.Range("N6").Value = 956127888.2484893
$ws.Range("O6").Value = 938536062.3753077
$ws.Range("P6").Value = 0.01800447008390172
$ws.Range("Q6").Value = 0.01834194406636469
